$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reference Lists")

# Fill in the missing "Data Set Status" dates for B3:B10 with the same
# date (12/5/2023, serial 45265) already used elsewhere in column B.
$ws.Range("B3:B10").Value = 45265

# Reuse the existing date formatting from B2 (numFmtId 14 / m/d/yyyy)
# rather than re-deriving a fresh custom number format.
$ws.Range("B2").Copy()
$ws.Range("B3:B10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active sheet's selection to match the new target range.
$ws.Activate()
$ws.Range("A11:A16").Select()
